$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 14; $r++) {
    foreach ($col in @("F", "G")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val.ToString().StartsWith("/mnt/z/")) {
            $cell.Value2 = $val.ToString().Replace("/mnt/z/", "/storeData/")
        }
    }
}
